# Update column G ("K") values on the active sheet (rows 2-37) to the
# newly-regenerated figures (K computed instead of the old "Strike#"-based
# values), per: "regen save_data to use K instead of Strike#, regen
# std/mean, calc and write s_vals".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 3
    4  = 5
    5  = 5
    6  = 5
    7  = 5
    8  = 7
    9  = 3
    10 = 3
    11 = 6
    12 = 8
    13 = 6
    14 = 7
    15 = 9
    16 = 3
    17 = 2
    18 = 1
    19 = 6
    20 = 5
    21 = 6
    22 = 9
    23 = 8
    24 = 14
    25 = 3
    26 = 5
    27 = 5
    28 = 6
    29 = 9
    30 = 6
    31 = 8
    32 = 6
    33 = 5
    34 = 4
    35 = 8
    36 = 3
    37 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
